# Insert a new "Protein Shakes" slide right before the existing
# "Feeding Tubes" slide (slide 3), pushing all following slides down
# by one position. We accomplish the insertion by duplicating slide 3
# (which already uses the desired "Two Content" layout with two
# half-width content placeholders), then overwriting the text on the
# original slide 3 with the new Protein Shakes content. The duplicate
# retains the original "Feeding Tubes" content and ends up at
# position 4, exactly where the old slide 3 content needs to be.

$p = $ppt.ActivePresentation

$original = $p.Slides.Item(3)
$duplicate = $original.Duplicate()

$original.Shapes.Item(1).TextFrame.TextRange.Text = "Protein Shakes"
$original.Shapes.Item(2).TextFrame.TextRange.Text = "Protein Shakes can provide protein with minimal sugar"
$original.Shapes.Item(3).TextFrame.TextRange.Text = "Portein Shakes "
